$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_A2 = @"
https://p16-images-sign-sg.tokopedia-static.net/tos-alisg-i-aphluv4xwc-sg/194dfcf7511141e4baa455e876099d65~tplv-aphluv4xwc-white-pad-v1:1600:1600.jpeg?lk3s=0ccea506&x-expires=1764166626&x-signature=SFftugZt8X28C8DZi6HWrg1iKUQ%3D&x-signature-webp=KVtQc8NB8Jo3bER3KzOxHfwhJMM%3D
"@
$ws.Range("A2").Value = $v_A2

$v_B2 = @"
Apple iPhone 15 Garansi Resmi - 128GB 256GB 512GB - Black, 128GB
"@
$ws.Range("B2").Value = $v_B2

$v_C2 = @"
Rp11.228.900
"@
$ws.Range("C2").Value = $v_C2

$v_D2 = @"
iSmile adalah APPLE ONLINE RESELLER yang khusus menjual produk APPLE ORIGINAL dan BERGARANSI RESMI.
iSmile adalah Tokopedia OFFICIAL STORE & REKOMENDASI SELLER APPLE PILIHAN Tokopedia
Customer Corporate / Project are welcome
iSmile menerima bulk order & bisa menerbitkan faktur pajak (konfirmasi dahulu sebelum transaksi)
Apakah produk Apple yang dijual iSmile sama seperti iBox, Digimap dan Apple Authorized Apple Reseller (AAR) lainnya?
Ya, Produk Apple yang kami jual ORIGINAL APPLE dan BERGARANSI RESMI INDONESIA, sama seperti yang dijual oleh iBox, Digimap, dll
Mengapa harga iSmile lebih mudah dibanding iBox, Digimap dan Apple Authorized Reseller (AAR) lainnya?
iSmile merupakan Apple Online Reseller yang TIDAK MEMILIKI TOKO FISIK, kami BEROPERASI secara EFISIEN karena tidak membayar sewa tempat yang mahal di mall, tidak menggaji banyak Sales Promotor. Sehingga HARGA yang kami tawarkan bisa LEBIH RENDAH.
Apakah perbedaan iSmile dengan Toko Online Lain (bukan Apple Authorized Reseller) di marketplace?
iSmile HANYA MENJUAL Produk APPLE ORIGINAL & BERGARANSI RESMI, Toko Online Lain (bukan Apple Authorized Reseller) banyak menjual Produk Apple Bergaransi International (Tidak Resmi), Refurbished, dll yang tidak bisa dijamin Originalitas-nya dan mungkin bermasalah pada saat kamu akan melakukan claim garansi di Apple Authorized Service Provider (AASP).
Bagaimana garansi Apple iSmile?
Semua produk Apple yang dijual iSmile BERGARANSI RESMI dan bisa diclaim di Apple Authorized Service Provider (AASP) Indonesia (iBox, Digimap, Mitracare, QCD Indonesia, Story-i, Hello) dan AASP Luar negeri.
"@
$ws.Range("D2").Value = $v_D2

$v_E2 = @"
Kota Administrasi Jakarta Pusat
"@
$ws.Range("E2").Value = $v_E2
